# Rename the "Test case Id" values in column A from the old C1..C8
# scheme to the new [LE]/01..[LE]/08 scheme (rows grouped per test case).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$renames = @{
    "A2"  = "[LE]/01"
    "A3"  = "[LE]/01"
    "A4"  = "[LE]/01"
    "A6"  = "[LE]/02"
    "A7"  = "[LE]/02"
    "A8"  = "[LE]/02"
    "A10" = "[LE]/03"
    "A11" = "[LE]/03"
    "A14" = "[LE]/04"
    "A15" = "[LE]/04"
    "A16" = "[LE]/04"
    "A18" = "[LE]/05"
    "A20" = "[LE]/06"
    "A22" = "[LE]/07"
    "A24" = "[LE]/08"
    "A25" = "[LE]/08"
}

foreach ($addr in $renames.Keys) {
    $ws.Range($addr).Value = $renames[$addr]
}
